$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 85532.34
$ws.Range("I15").Value = 85532.34
$ws.Range("K15").Value = 256597.02
$ws.Range("M15").Value = -256428.02
$ws.Range("H19").Value = 636.9286
$ws.Range("I19").Value = 481.3846
$ws.Range("J19").Value = 771.73334
$ws.Range("K19").Value = 481.3846
$ws.Range("L19").Value = 771.73334
$ws.Range("M19").Value = -306.3846
$ws.Range("N19").Value = -1121.73334
$ws.Range("H28").Value = 371216.03
$ws.Range("I28").Value = 529601.1
$ws.Range("J28").Value = 1650.7778
$ws.Range("K28").Value = 529601.1
$ws.Range("L28").Value = 1650.7778
$ws.Range("M28").Value = -529116.1
$ws.Range("N28").Value = -2620.7778
$ws.Range("H38").Value = 97.85714
$ws.Range("J38").Value = 490
$ws.Range("L38").Value = 1470
$ws.Range("N38").Value = -2214
$ws.Range("H58").Value = 3981.1365
$ws.Range("I58").Value = 98.5
$ws.Range("J58").Value = 7216.6665
$ws.Range("K58").Value = 295.5
$ws.Range("L58").Value = 21649.9995
$ws.Range("M58").Value = -145.5
$ws.Range("N58").Value = -21949.9995
$ws.Range("H107").Value = 505591.4
$ws.Range("I107").Value = 529633.9
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 529633.9
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = -527713.9
$ws.Range("N107").Value = -4540
$ws.Range("H123").Value = 65600
$ws.Range("J123").Value = 65600
$ws.Range("L123").Value = 65600
$ws.Range("N123").Value = -75400
$ws.Range("H124").Value = 29000
$ws.Range("J124").Value = 29000
$ws.Range("L124").Value = 29000
$ws.Range("N124").Value = -38820
$ws.Range("H129").Value = 1227.1428
$ws.Range("I129").Value = 390.66666
$ws.Range("J129").Value = 1854.5
$ws.Range("K129").Value = 1171.99998
$ws.Range("L129").Value = 5563.5
$ws.Range("M129").Value = 3828.00002
$ws.Range("N129").Value = -15563.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H97").Value = 8339.923000000001
$ws.Range("I97").Value = 10591.9
$ws.Range("J97").Value = 833.3333
$ws.Range("K97").Value = 10591.9
$ws.Range("L97").Value = 833.3333
$ws.Range("M97").Value = -10095.9
$ws.Range("N97").Value = -1825.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20240.834
$ws.Range("I82").Value = 7861.25
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 7861.25
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -7478.25
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 20240.834
$ws.Range("I85").Value = 7861.25
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 7861.25
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -6535.25
$ws.Range("N85").Value = -47652

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 334000
$ws.Range("J4").Value = 334000
$ws.Range("L4").Value = 334000
$ws.Range("N4").Value = -334224

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5666.6665
$ws.Range("I56").Value = 5666.6665
$ws.Range("K56").Value = 5666.6665
$ws.Range("M56").Value = -5136.6665
$ws.Range("H94").Value = 2005.75
$ws.Range("I94").Value = 1007.6667
$ws.Range("K94").Value = 3023.0001
$ws.Range("M94").Value = -2347.0001
$ws.Range("H129").Value = 2120.5557
$ws.Range("I129").Value = 3676.3333
$ws.Range("J129").Value = 1342.6666
$ws.Range("K129").Value = 11028.9999
$ws.Range("L129").Value = 4027.9998
$ws.Range("M129").Value = -6028.999899999999
$ws.Range("N129").Value = -14027.9998
$ws.Range("H131").Value = 8773773
$ws.Range("I131").Value = 943.3333
$ws.Range("J131").Value = 9525730
$ws.Range("K131").Value = 2829.9999
$ws.Range("L131").Value = 28577190
$ws.Range("M131").Value = 2210.0001
$ws.Range("N131").Value = -28587270

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 40002816
$ws.Range("I80").Value = 2656.0625
$ws.Range("J80").Value = 111114210
$ws.Range("K80").Value = 2656.0625
$ws.Range("L80").Value = 111114210
$ws.Range("M80").Value = -1658.0625
$ws.Range("N80").Value = -111116206
$ws.Range("H83").Value = 40002816
$ws.Range("I83").Value = 2656.0625
$ws.Range("J83").Value = 111114210
$ws.Range("K83").Value = 13280.3125
$ws.Range("L83").Value = 555571050
$ws.Range("M83").Value = -8288.3125
$ws.Range("N83").Value = -555581034
$ws.Range("H122").Value = 2728.9285
$ws.Range("J122").Value = 2266.6667
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("N122").Value = -11700.0001
$ws.Range("H123").Value = 11462.608
$ws.Range("J123").Value = 11462.608
$ws.Range("L123").Value = 11462.608
$ws.Range("N123").Value = -16362.608

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3400
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452
$ws.Range("H22").Value = 634.1429000000001
$ws.Range("I22").Value = 659.8333
$ws.Range("J22").Value = 480
$ws.Range("K22").Value = 659.8333
$ws.Range("L22").Value = 480
$ws.Range("M22").Value = -364.8333
$ws.Range("N22").Value = -1070
$ws.Range("H24").Value = 4000
$ws.Range("J24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("N24").Value = -4686
$ws.Range("H27").Value = 634.1429000000001
$ws.Range("I27").Value = 659.8333
$ws.Range("J27").Value = 480
$ws.Range("K27").Value = 659.8333
$ws.Range("L27").Value = 480
$ws.Range("M27").Value = -552.8333
$ws.Range("N27").Value = -694
$ws.Range("H40").Value = 3388.8462
$ws.Range("I40").Value = 2777.5
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 2777.5
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -2641.5
$ws.Range("N40").Value = -3772
$ws.Range("H61").Value = 1270.2858
$ws.Range("I61").Value = 724.25
$ws.Range("J61").Value = 1998.3334
$ws.Range("K61").Value = 724.25
$ws.Range("L61").Value = 1998.3334
$ws.Range("M61").Value = -522.25
$ws.Range("N61").Value = -2402.3334
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H113").Value = 1270.2858
$ws.Range("I113").Value = 724.25
$ws.Range("J113").Value = 1998.3334
$ws.Range("K113").Value = 724.25
$ws.Range("L113").Value = 1998.3334
$ws.Range("M113").Value = 1445.75
$ws.Range("N113").Value = -6338.3334
$ws.Range("H132").Value = 3471.8462
$ws.Range("I132").Value = 2173.9443
$ws.Range("J132").Value = 6392.125
$ws.Range("K132").Value = 6521.8329
$ws.Range("L132").Value = 19176.375
$ws.Range("M132").Value = -3991.8329
$ws.Range("N132").Value = -24236.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 33500
$ws.Range("I9").Value = 33500
$ws.Range("K9").Value = 33500
$ws.Range("M9").Value = -33360
$ws.Range("H12").Value = 11860
$ws.Range("I12").Value = 50000
$ws.Range("J12").Value = 2325
$ws.Range("K12").Value = 50000
$ws.Range("L12").Value = 2325
$ws.Range("M12").Value = -49858
$ws.Range("N12").Value = -2609
